$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: fertilization-event-2 columns (S:V) become explicit empty strings ---
$rng = $ws.Range("S2:V2")
$rng.Value = "'"
$rng.Style = "Normal"

# --- Rows 3-20: add a second fertilization event (S:V) using the new 2024-04-04 date ---
$rng = $ws.Range("S3:S20")
$rng.Value = "'2024-04-04"
$rng.Style = "Normal"

$rng = $ws.Range("T3:T20")
$rng.Value = "2-2-2"
$rng.Style = "Normal"

$rng = $ws.Range("U3:U20")
$rng.Value = "Espoma- Organic Grow!"
$rng.Style = "Normal"

$rng = $ws.Range("V3:V20")
$rng.Value = "half"
$rng.Style = "Normal"

# --- Rows 21-33: fertilization-event-1 columns (O:R) now populated with the new date ---
$rng = $ws.Range("O21:O33")
$rng.Value = "'2024-04-04"
$rng.Style = "Normal"

$rng = $ws.Range("P21:P33")
$rng.Value = "2-2-2"
$rng.Style = "Normal"

$rng = $ws.Range("Q21:Q33")
$rng.Value = "Espoma- Organic Grow!"
$rng.Style = "Normal"

$rng = $ws.Range("R21:R33")
$rng.Value = "half"
$rng.Style = "Normal"

# --- Rows 21-33: fertilization-event-2 columns (S:V) become explicit empty strings ---
$rng = $ws.Range("S21:V33")
$rng.Value = "'"
$rng.Style = "Normal"
